$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# The "Periodo Mora" data table (rows 16-22) is being re-sorted so the
# periods read in ascending order (2405 .. 2411) instead of descending
# (2411 .. 2405). The "Valor Mora" (column F) value travels with its
# period: period 2411 carries 24266, all other periods carry 52000.

$periods = @("2405", "2406", "2407", "2408", "2409", "2410", "2411")
$valores = @(52000, 52000, 52000, 52000, 52000, 52000, 24266)

for ($i = 0; $i -lt $periods.Length; $i++) {
    $row = 16 + $i
    $ws.Cells.Item($row, 5).Value = $periods[$i]
    $ws.Cells.Item($row, 6).Value = $valores[$i]
}
